# Server.xlsx / "Property" sheet update
# - Replace the placeholder "127.0.0.1" IP values (column F, rows 2-6) with
#   distinct real addresses.
# - Remove the extra "GameServer_2" row (row 7), leaving the row blank but
#   keeping its original cell formatting.
# - Resize column F (IP) to fit its new, wider content.
# - Leave the selection on row 7, matching the state after the row contents
#   were cleared in the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP column with the new addresses
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Drop the duplicate "GameServer_2" row; clear contents but keep formatting
$ws.Range("A7:H7").ClearContents()

# Column F needs to be widened to fit "192.168.1.11x"
$ws.Columns("F").ColumnWidth = 14.285714285714286

# Match the editor's final selection (entire row 7)
$ws.Rows(7).Select()
